$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1 (header row): rename B1:G1 to reflect the reshuffled columns.
# (J1's shared-string index also shifts in the diff, but the text itself
#  stays "surplus_protein" - nothing to change there.)
# ---------------------------------------------------------------------------
$ws.Range("B1").Value = "get_up"
$ws.Range("C1").Value = "bedtime"
$ws.Range("D1").Value = "activity"
$ws.Range("E1").Value = "stress_level"
$ws.Range("F1").Value = "headpain"
$ws.Range("G1").Value = "Painkiller"

# ---------------------------------------------------------------------------
# Row 2
#   before: B2=none(date-fmt) C2=yes(date-fmt) D2=<7(text) E2=<21(text) F2=uni(text) G2=5(text,number)
#   after:  B2=<7(text) C2=<21(text) D2=uni(text) E2=5(text,number) F2=none(date-fmt) G2=yes(date-fmt)
# Use stable, untouched helper cells as format sources:
#   H2 always keeps its original text ("@") style -> donor for text-style cells
#   A2 always keeps its original date style       -> donor for date-style cells
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy()
$ws.Range("B2:C2").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("F2:G2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B2").Value = "<7"
$ws.Range("C2").Value = "<21"
$ws.Range("D2").Value = "uni"
$ws.Range("F2").Value = "none"
$ws.Range("G2").Value = "yes"

# E2 must hold a genuine number (5) while styled as text ("@"). Writing a
# number straight into an "@"-formatted cell gets stored as text, so clear
# the format first, write the number, then re-apply the text format last.
$ws.Range("E2").ClearFormats()
$ws.Range("E2").Value = 5
$ws.Range("H2").Copy()
$ws.Range("E2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Row 3
#   before: B3=little(date-fmt)           D3=7-9 E3=21-23 F3=lab G3=0-10 (all text-fmt)
#   after:  B3=7-9 C3=21-23 D3=lab E3=0-10 (text-fmt)   F3=little(date-fmt)   G3 removed
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy()
$ws.Range("B3:C3").PasteSpecial(-4122)
$ws.Range("A2").Copy()
$ws.Range("F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B3").Value = "7-9"
$ws.Range("C3").Value = "21-23"
$ws.Range("D3").Value = "lab"
$ws.Range("E3").Value = "0-10"
$ws.Range("F3").Value = "little"
$ws.Range("G3").Clear()

# ---------------------------------------------------------------------------
# Row 4
#   before: B4=medium(no style)   D4=>9 E4=>23 F4=learning (text-fmt)   G4=(empty, text-fmt)
#   after:  B4=>9 C4=>23 D4=learning (text-fmt)   E4=(empty, text-fmt)   F4=medium(no style)   G4 removed
# ---------------------------------------------------------------------------
$ws.Range("H2").Copy()
$ws.Range("B4:C4").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F4").ClearFormats()

$ws.Range("B4").Value = ">9"
$ws.Range("C4").Value = ">23"
$ws.Range("D4").Value = "learning"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "medium"
$ws.Range("G4").Clear()

# ---------------------------------------------------------------------------
# Row 5
#   before: B5=massive(no style)   D5=(empty) E5=(empty) F5=housework (text-fmt)   G5=(empty, text-fmt)
#   after:  D5=housework (text-fmt)   E5=(empty, text-fmt)   F5=massive(no style)   B5,G5 removed
# ---------------------------------------------------------------------------
$ws.Range("F5").ClearFormats()

$ws.Range("B5").Clear()
$ws.Range("D5").Value = "housework"
$ws.Range("F5").Value = "massive"
$ws.Range("G5").Clear()

# ---------------------------------------------------------------------------
# Row 6
#   before: D6=(empty) E6=(empty) F6=sparetime (text-fmt)   G6=(empty, text-fmt)
#   after:  D6=sparetime (text-fmt)   E6=(empty, text-fmt)   F6,G6 removed
# ---------------------------------------------------------------------------
$ws.Range("D6").Value = "sparetime"
$ws.Range("F6").Clear()
$ws.Range("G6").Clear()

# ---------------------------------------------------------------------------
# Row 9: drop the leftover "TODO:" / "X" / "Zahlen einfuegen" helper notes.
#   before: A9="TODO:" C9="X" G9="Zahlen einfuegen"(text-fmt)
#   after:  A9,C9 removed;  G9 value cleared but keeps its text-fmt style
# ---------------------------------------------------------------------------
$ws.Range("A9").Clear()
$ws.Range("C9").Clear()
$ws.Range("G9").ClearContents()

# ---------------------------------------------------------------------------
# Rows 10-15: drop now-unused empty placeholder cells.
# ---------------------------------------------------------------------------
$ws.Range("D10").Clear()
$ws.Range("E10").Clear()

$ws.Range("D11:G11").Clear()
$ws.Range("D12:G12").Clear()
$ws.Range("D13:G13").Clear()
$ws.Range("D14:G14").Clear()

$ws.Range("F15").Clear()
$ws.Range("G15").Clear()

# ---------------------------------------------------------------------------
# Selection moves to F12.
# ---------------------------------------------------------------------------
$ws.Range("F12").Select()
